$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1.5888322343738148
$ws.Range("F2").Value = 0.014933622914248774
$ws.Range("G2").Value = 0.14276476253786627
$ws.Range("H2").Value = 0.12523224784023357

$ws.Range("E3").Value = 1.5358711598946877
$ws.Range("F3").Value = 0.023794127089644377
$ws.Range("G3").Value = 0.14276476253786627
$ws.Range("H3").Value = 0.12523224784023357

$ws.Range("E4").Value = -1.3438956514398925
$ws.Range("F4").Value = 0.10056925996204934
$ws.Range("G4").Value = 0.3415559772296015
$ws.Range("H4").Value = 0.29961050634175573

$ws.Range("E5").Value = -1.3302944201672433
$ws.Range("F5").Value = 0.11385199240986717
$ws.Range("G5").Value = 0.3415559772296015
$ws.Range("H5").Value = 0.29961050634175573
